$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.236.84"
$ws.Range("E2").Value = "  -3.44%  "

$ws.Range("D3").Value = "'3.140.36"
$ws.Range("E3").Value = "  -2.36%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'598.72"
$ws.Range("E5").Value = "  -1.45%  "

$ws.Range("D6").Value = "'150.31"
$ws.Range("E6").Value = "  -4.96%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "'3.141.44"
$ws.Range("E8").Value = "  -2.34%  "

$ws.Range("D9").Value = "'0.535"
$ws.Range("E9").Value = "  -2.74%  "

$ws.Range("D10").Value = "'0.154"
$ws.Range("E10").Value = "  -4.24%  "

$ws.Range("E11").Value = "  -0.86%  "

$ws.Range("D12").Value = "'0.481"
$ws.Range("E12").Value = "  -4.29%  "

$ws.Range("D13").Value = "'0.0000261"
$ws.Range("E13").Value = "  -3.25%  "

$ws.Range("D14").Value = "'37.09"
$ws.Range("E14").Value = "  -4.10%  "

$ws.Range("D15").Value = "'3.614.58"
$ws.Range("E15").Value = "  -3.47%  "

$ws.Range("D16").Value = "'64.344.30"
$ws.Range("E16").Value = "  -3.43%  "

$ws.Range("E17").Value = "  +0.45%  "

$ws.Range("D18").Value = "'3.134.45"
$ws.Range("E18").Value = "  -2.68%  "

$ws.Range("D19").Value = "'7.06"
$ws.Range("E19").Value = "  -4.00%  "

$ws.Range("D20").Value = "'484.25"
$ws.Range("E20").Value = "  -4.40%  "

$ws.Range("D21").Value = "'14.79"
$ws.Range("E21").Value = "  -2.33%  "

$ws.Range("D22").Value = "'0.715"
$ws.Range("E22").Value = "  -2.28%  "

$ws.Range("D23").Value = "'7.84"
$ws.Range("E23").Value = "  -2.03%  "

$ws.Range("D24").Value = "'13.93"
$ws.Range("E24").Value = "  -4.41%  "

$ws.Range("D25").Value = "'85.17"
$ws.Range("E25").Value = "  +0.40%  "

$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("E27").Value = "  -2.36%  "

$ws.Range("D28").Value = "'8.70"
$ws.Range("E28").Value = "  -4.43%  "

$ws.Range("D29").Value = "'2.27"
$ws.Range("E29").Value = "  -3.91%  "

$ws.Range("E30").Value = "  +2.31%  "

$ws.Range("D31").Value = "'7.12"
$ws.Range("E31").Value = "  +1.69%  "

$ws.Range("E32").Value = "  -7.21%  "

$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("D34").Value = "'26.91"
$ws.Range("E34").Value = "  -4.21%  "

$ws.Range("E35").Value = "  -5.96%  "

$ws.Range("D36").Value = "'6.13"
$ws.Range("E36").Value = "  -5.27%  "

$ws.Range("E37").Value = "  +8.31%  "

$ws.Range("D38").Value = "'54.71"
$ws.Range("E38").Value = "  -1.22%  "

$ws.Range("D39").Value = "'0.0₃0753"
$ws.Range("E39").Value = "  -2.19%  "

$ws.Range("D40").Value = "'451.72"
$ws.Range("E40").Value = "  -9.90%  "

$ws.Range("E41").Value = "  -4.61%  "

$ws.Range("D42").Value = "'0.0404"
$ws.Range("E42").Value = "  -3.98%  "

$ws.Range("D43").Value = "'8.56"
$ws.Range("E43").Value = "  -1.76%  "

$ws.Range("D44").Value = "'2.43"
$ws.Range("E44").Value = "  -0.74%  "

$ws.Range("D45").Value = "'2.899.69"
$ws.Range("E45").Value = "  -0.39%  "

$ws.Range("D46").Value = "'0.276"
$ws.Range("E46").Value = "  -6.93%  "

$ws.Range("D47").Value = "'27.08"
$ws.Range("E47").Value = "  -3.80%  "

$ws.Range("D48").Value = "'0.999"
$ws.Range("E48").Value = "  +0.02%  "

$ws.Range("E49").Value = "  +0.32%  "

$ws.Range("D50").Value = "'2.34"
$ws.Range("E50").Value = "  -2.84%  "

$ws.Range("D51").Value = "'2.51"
$ws.Range("E51").Value = "  +0.31%  "
